$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosCP")

# ---------------------------------------------------------------------------
# Row 5 used to hold the lone "CP004" test case (only column A populated).
# It now becomes a full test case row, just like rows 2-4, with a new
# search/login test ("CP004_buscarPorNombre") plus an extra expected-result
# column (E) holding "Venta cancelada".
# ---------------------------------------------------------------------------

# Give A5:D5 the same formatting rows 2-4 already use (font/border/fill)
# before writing the new values, so the row keeps a consistent look.
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)   # xlPasteFormats

# Fill in the new test-case data. E5 is written before D5 so that the
# shared-string table ends up with "Venta cancelada" interned ahead of
# "rodolfo", matching how the workbook was actually edited.
$ws.Range("A5").Value = "CP004_buscarPorNombre"
$ws.Range("B5").Value = "spawnro2@gmail.com"
$ws.Range("C5").Value = "diego170915#"
$ws.Range("E5").Value = "Venta cancelada"
$ws.Range("D5").Value = "rodolfo"

# E5 gets its own look: the same green monospace font used for the other
# data columns, combined with a thin box border and vertically centered
# text (like the boxed cells in row 6).
$ws.Range("C2").Copy()
$ws.Range("E5").PasteSpecial(-4122)      # pulls in font + vertical-center alignment
$e5 = $ws.Range("E5")
$e5.Borders.LineStyle = 1
$e5.Borders.Weight = 2

# Column E now holds real text, so let it auto-size like the other
# best-fit columns (A-D).
$ws.Columns.Item(5).AutoFit()

# Restore the normal paste-mode marching ants / clipboard state.
$excel.CutCopyMode = $false

# Leave the selection where the author's session ended up.
[void]$ws.Range("D12").Select()
